$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per-row: D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$cols = @(4, 10, 11, 12, 13, 16)

# Snapshot the current (pre-edit) values for every data row (2..25) before
# any writes happen, keyed by row number, so the row-to-row shuffle below
# reads consistent "before" data regardless of write order.
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Target row -> source row mapping (row 14 keeps its original data).
$mapping = @{
    2  = 21
    3  = 13
    4  = 3
    5  = 6
    6  = 20
    7  = 8
    8  = 22
    9  = 16
    10 = 11
    11 = 24
    12 = 10
    13 = 25
    15 = 12
    16 = 5
    17 = 4
    18 = 19
    19 = 9
    20 = 2
    21 = 17
    22 = 18
    23 = 7
    24 = 23
    25 = 15
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $sourceVals[$c]
    }
}
